# The "find and replace" worksheet holds a B3:H20 grid of "A"/"Z" shared-
# string labels used for a find/replace exercise. This edit blanks out a
# scattered set of cells (turning them from an "A" label into a truly
# empty cell) so the exercise has gaps to fill in / find against.
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("find and replace")

$cellsToClear = @(
    "D3", "F3",
    "B4", "E4", "G4",
    "D5", "F5",
    "C6", "G6",
    "B7", "E7", "H7",
    "D8",
    "B9", "F9", "H9",
    "C10", "D10", "E10", "F10", "H10",
    "B11", "C11",
    "E12",
    "C13",
    "F14", "G14",
    "B15", "H15",
    "F16", "G16",
    "B17", "D17",
    "F18", "H18",
    "C19", "D19", "E19", "F19",
    "B20", "D20", "F20"
)

foreach ($ref in $cellsToClear) {
    $ws.Range($ref).ClearContents()
}
